$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 230.84616
$ws.Range("I5").Value = 51.375
$ws.Range("J5").Value = 518
$ws.Range("K5").Value = 51.375
$ws.Range("L5").Value = 518
$ws.Range("M5").Value = 63.625
$ws.Range("N5").Value = -748
$ws.Range("H112").Value = 39359.617
$ws.Range("I112").Value = 430
$ws.Range("J112").Value = 48628.57
$ws.Range("K112").Value = 1290
$ws.Range("L112").Value = 145885.71
$ws.Range("M112").Value = -182
$ws.Range("N112").Value = -148101.71
$ws.Range("H113").Value = 2107.2432
$ws.Range("I113").Value = 2172.456
$ws.Range("J113").Value = 1888.5883
$ws.Range("K113").Value = 2172.456
$ws.Range("L113").Value = 1888.5883
$ws.Range("M113").Value = 1081.544
$ws.Range("N113").Value = -8396.588299999999
$ws.Range("H129").Value = 399.9
$ws.Range("I129").Value = 299.8889
$ws.Range("K129").Value = 899.6667
$ws.Range("M129").Value = 4100.3333
$ws.Range("H137").Value = 4106.231
$ws.Range("I137").Value = 3776.476
$ws.Range("J137").Value = 5491.2
$ws.Range("K137").Value = 11329.428
$ws.Range("L137").Value = 16473.6
$ws.Range("M137").Value = -8779.428
$ws.Range("N137").Value = -21573.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 941.6667
$ws.Range("I2").Value = 816.6667
$ws.Range("J2").Value = 1066.6666
$ws.Range("K2").Value = 816.6667
$ws.Range("L2").Value = 1066.6666
$ws.Range("M2").Value = -703.6667
$ws.Range("N2").Value = -1292.6666
$ws.Range("H32").Value = 2805.6191
$ws.Range("I32").Value = 1786.8
$ws.Range("J32").Value = 7899.7144
$ws.Range("K32").Value = 1786.8
$ws.Range("L32").Value = 7899.7144
$ws.Range("M32").Value = -1499.8
$ws.Range("N32").Value = -8473.714400000001
$ws.Range("H61").Value = 8290.666999999999
$ws.Range("I61").Value = 1795.8462
$ws.Range("J61").Value = 50507
$ws.Range("K61").Value = 1795.8462
$ws.Range("L61").Value = 50507
$ws.Range("M61").Value = -1583.8462
$ws.Range("N61").Value = -50931
$ws.Range("H74").Value = 1414.7059
$ws.Range("I74").Value = 603.1429000000001
$ws.Range("J74").Value = 1982.8
$ws.Range("K74").Value = 603.1429000000001
$ws.Range("L74").Value = 1982.8
$ws.Range("M74").Value = 270.8570999999999
$ws.Range("N74").Value = -3730.8
$ws.Range("H77").Value = 1414.7059
$ws.Range("I77").Value = 603.1429000000001
$ws.Range("J77").Value = 1982.8
$ws.Range("K77").Value = 3015.7145
$ws.Range("L77").Value = 9914
$ws.Range("M77").Value = 1352.2855
$ws.Range("N77").Value = -18650
$ws.Range("H116").Value = 941.6667
$ws.Range("I116").Value = 816.6667
$ws.Range("J116").Value = 1066.6666
$ws.Range("K116").Value = 816.6667
$ws.Range("L116").Value = 1066.6666
$ws.Range("M116").Value = 1477.3333
$ws.Range("N116").Value = -5654.6666
$ws.Range("H132").Value = 2232.6736
$ws.Range("I132").Value = 2042.575
$ws.Range("J132").Value = 3077.5557
$ws.Range("K132").Value = 6127.725
$ws.Range("L132").Value = 9232.667099999999
$ws.Range("M132").Value = -3597.725
$ws.Range("N132").Value = -14292.6671
$ws.Range("H136").Value = 8290.666999999999
$ws.Range("I136").Value = 1795.8462
$ws.Range("J136").Value = 50507
$ws.Range("K136").Value = 5387.5386
$ws.Range("L136").Value = 151521
$ws.Range("M136").Value = -2837.5386
$ws.Range("N136").Value = -156621

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 941.6667
$ws.Range("I3").Value = 816.6667
$ws.Range("J3").Value = 1066.6666
$ws.Range("K3").Value = 816.6667
$ws.Range("L3").Value = 1066.6666
$ws.Range("M3").Value = -702.6667
$ws.Range("N3").Value = -1294.6666
$ws.Range("H94").Value = 1365.75
$ws.Range("I94").Value = 1347.9412
$ws.Range("J94").Value = 1466.6666
$ws.Range("K94").Value = 1347.9412
$ws.Range("L94").Value = 1466.6666
$ws.Range("M94").Value = -896.9412
$ws.Range("N94").Value = -2368.6666
$ws.Range("H107").Value = 1131.6
$ws.Range("I107").Value = 1136.4117
$ws.Range("J107").Value = 1104.3334
$ws.Range("K107").Value = 1136.4117
$ws.Range("L107").Value = 1104.3334
$ws.Range("M107").Value = 783.5882999999999
$ws.Range("N107").Value = -4944.3334
$ws.Range("H113").Value = 4680
$ws.Range("I113").Value = 4680
$ws.Range("K113").Value = 4680
$ws.Range("M113").Value = -2510
$ws.Range("H134").Value = 2364.4546
$ws.Range("I134").Value = 1600.9
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 4802.700000000001
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -2267.700000000001
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1224.75
$ws.Range("I16").Value = 1224.75
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1224.75
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 34432.35
$ws.Range("I31").Value = 48967.22
$ws.Range("J31").Value = 4041.2727
$ws.Range("K31").Value = 48967.22
$ws.Range("L31").Value = 4041.2727
$ws.Range("M31").Value = -48672.22
$ws.Range("N31").Value = -4631.2727
$ws.Range("H34").Value = 34432.35
$ws.Range("I34").Value = 48967.22
$ws.Range("J34").Value = 4041.2727
$ws.Range("K34").Value = 48967.22
$ws.Range("L34").Value = 4041.2727
$ws.Range("M34").Value = -48765.22
$ws.Range("N34").Value = -4445.2727
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").ClearContents()
$ws.Range("H106").Value = 150000
$ws.Range("J106").Value = 150000
$ws.Range("L106").Value = 150000
$ws.Range("N106").Value = -152524
$ws.Range("H107").Value = 1582.2222
$ws.Range("I107").Value = 1587.0588
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1587.0588
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 332.9412
$ws.Range("N107").Value = -5340
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H110").Value = 50702
$ws.Range("J110").Value = 50702
$ws.Range("L110").Value = 50702
$ws.Range("N110").Value = -58882
$ws.Range("H111").Value = 50702
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("H113").Value = 1224.75
$ws.Range("I113").Value = 1224.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1224.75
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 7102.95
$ws.Range("I132").Value = 1488.2307
$ws.Range("J132").Value = 17530.285
$ws.Range("K132").Value = 4464.6921
$ws.Range("L132").Value = 52590.855
$ws.Range("M132").Value = -1934.6921
$ws.Range("N132").Value = -57650.855
$ws.Range("H134").Value = 8612.526
$ws.Range("I134").Value = 2538.5454
$ws.Range("J134").Value = 16964.25
$ws.Range("K134").Value = 7615.6362
$ws.Range("L134").Value = 50892.75
$ws.Range("M134").Value = -5080.6362
$ws.Range("N134").Value = -55962.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5190.5264
$ws.Range("I70").Value = 5288
$ws.Range("J70").Value = 4825
$ws.Range("K70").Value = 5288
$ws.Range("L70").Value = 4825
$ws.Range("M70").Value = -5018
$ws.Range("N70").Value = -5365
$ws.Range("H73").Value = 5190.5264
$ws.Range("I73").Value = 5288
$ws.Range("J73").Value = 4825
$ws.Range("K73").Value = 5288
$ws.Range("L73").Value = 4825
$ws.Range("M73").Value = -4352
$ws.Range("N73").Value = -6697
$ws.Range("H113").Value = 1003.3077
$ws.Range("I113").Value = 822.1667
$ws.Range("J113").Value = 1158.5714
$ws.Range("K113").Value = 822.1667
$ws.Range("L113").Value = 1158.5714
$ws.Range("M113").Value = 1347.8333
$ws.Range("N113").Value = -5498.5714
$ws.Range("H122").Value = 1569.4615
$ws.Range("I122").Value = 1350.25
$ws.Range("J122").Value = 2300.1667
$ws.Range("K122").Value = 4050.75
$ws.Range("L122").Value = 6900.500100000001
$ws.Range("M122").Value = -1600.75
$ws.Range("N122").Value = -11800.5001
$ws.Range("H132").Value = 2081.7585
$ws.Range("I132").Value = 1651.4348
$ws.Range("J132").Value = 3731.3333
$ws.Range("K132").Value = 4954.3044
$ws.Range("L132").Value = 11193.9999
$ws.Range("M132").Value = -2424.3044
$ws.Range("N132").Value = -16253.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2564.3572
$ws.Range("I40").Value = 2325.0833
$ws.Range("K40").Value = 2325.0833
$ws.Range("M40").Value = -2189.0833
$ws.Range("H132").Value = 7498.026
$ws.Range("I132").Value = 6763.421
$ws.Range("J132").Value = 8195.9
$ws.Range("K132").Value = 20290.263
$ws.Range("L132").Value = 24587.7
$ws.Range("M132").Value = -17760.263
$ws.Range("N132").Value = -29647.7
$ws.Range("H136").Value = 6421.1
$ws.Range("I136").Value = 5072.241
$ws.Range("J136").Value = 8283.809999999999
$ws.Range("K136").Value = 15216.723
$ws.Range("L136").Value = 24851.43
$ws.Range("M136").Value = -12666.723
$ws.Range("N136").Value = -29951.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6474.0527
$ws.Range("I132").Value = 515.6923
$ws.Range("J132").Value = 19383.834
$ws.Range("K132").Value = 1547.0769
$ws.Range("L132").Value = 58151.50199999999
$ws.Range("M132").Value = 982.9231
$ws.Range("N132").Value = -63211.50199999999
